$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false,
                             $true, 1, $false, $new, 2)
}

Replace-Text "782×8=" "525×7="
Replace-Text "542×4=" "924×7="
Replace-Text "775×9=" "334×5="
Replace-Text "988×8=" "429×4="
Replace-Text "543×2=" "399×6="
Replace-Text "292×9=" "153×5="
Replace-Text "489×6=" "767×8="
Replace-Text "475×9=" "857×6="
Replace-Text "935×9=" "885×7="
Replace-Text "605×7=" "239×8="
Replace-Text "861×9=" "272×2="
Replace-Text "538×4=" "574×3="
Replace-Text "627×2=" "972×6="
Replace-Text "472×5=" "257×5="
Replace-Text "919×4=" "742×9="
Replace-Text "849×7=" "418×6="
Replace-Text "278×7=" "812×9="
Replace-Text "509×3=" "946×8="
Replace-Text "610×3=" "990×6="
Replace-Text "178×7=" "863×4="
Replace-Text "536×9=" "927×7="
Replace-Text "560×6=" "160×3="
Replace-Text "636×4=" "419×3="
Replace-Text "266×2=" "642×9="
Replace-Text "656×9=" "245×4="
